$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 - Bitcoin
$ws.Range("D2").Value = "25.920.62"
$ws.Range("E2").Value = "  -0.12%  "

# Row 3 - Ethereum
$ws.Range("D3").Value = "1.632.89"
$ws.Range("E3").Value = "  -0.30%  "

# Row 4 - TetherUSD
$ws.Range("E4").Value = "  +0.06%  "

# Row 5 - BNB
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.00"
$ws.Range("E5").Value = "  +0.73%  "

# Row 6 - XRP
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5106"
$ws.Range("E6").Value = "  +0.14%  "

# Row 7 - USDC
$ws.Range("E7").Value = "  +0.13%  "

# Row 8 - Cardano
$ws.Range("E8").Value = "  +0.65%  "

# Row 9 - Dogecoin
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06344"
$ws.Range("E9").Value = "  -0.25%  "

# Row 10 - Solana
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "19.48"
$ws.Range("E10").Value = "  -0.04%  "

# Row 11 - TRON
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07782"
$ws.Range("E11").Value = "  +0.45%  "

# Row 12 - Polkadot
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "4.276"
$ws.Range("E12").Value = "  -0.04%  "

# Row 13 - WrappedEther
$ws.Range("D13").Value = "1.629.21"
$ws.Range("E13").Value = "  -0.81%  "

# Row 14 - WrappedliquidstakedEther2.0
$ws.Range("D14").Value = "1.857.99"
$ws.Range("E14").Value = "  -0.37%  "

# Row 15 - Polygon
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5510"
$ws.Range("E15").Value = "  +1.50%  "

# Row 16 - Litecoin
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.93"
$ws.Range("E16").Value = "  -0.51%  "

# Row 17 - ShibaInu
$ws.Range("D17").Value = "0.0₅7653"
$ws.Range("E17").Value = "  -1.36%  "

# Row 18 - WrappedBTC
$ws.Range("D18").Value = "25.945.34"
$ws.Range("E18").Value = "  -0.03%  "

# Row 19 - Dai
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.002"
$ws.Range("E19").Value = "  +0.17%  "

# Row 20 - BitcoinCash
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "194.83"
$ws.Range("E20").Value = "  -0.71%  "

# Row 21 - Uniswap
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "4.415"
$ws.Range("E21").Value = "  -0.05%  "

# Row 22 - Avalanche
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.865"
$ws.Range("E22").Value = "  -0.50%  "

# Row 23 - Chainlink
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.049"
$ws.Range("E23").Value = "  +0.27%  "

# Row 24 - BinanceUSD
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.003"
$ws.Range("E24").Value = "  +0.02%  "

# Row 25 - Toncoin
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.891"
$ws.Range("E25").Value = "  +0.87%  "

# Row 26 - Monero
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "142.09"
$ws.Range("E26").Value = "  +0.83%  "

# Row 27 - Stellar
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "0.1258"
$ws.Range("E27").Value = "  +5.52%  "

# Row 28 - Cosmos -> EthereumClassic
$ws.Range("B28").Value = "EthereumClassic"
$ws.Range("C28").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "15.59"
$ws.Range("E28").Value = "  -0.12%  "

# Row 29 - EthereumClassic -> Cosmos
$ws.Range("B29").Value = "Cosmos"
$ws.Range("C29").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "6.749"
$ws.Range("E29").Value = "  -1.34%  "

# Row 30 - PancakeSwap
$ws.Range("E30").Value = "  +0.71%  "

# Row 31 - Hedera
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.04887"
$ws.Range("E31").Value = "  -0.95%  "

# Row 32 - InternetComputer(DFINITY)
$ws.Range("E32").Value = "  -0.08%  "

# Row 33 - Filecoin
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.190"
$ws.Range("E33").Value = "  +0.47%  "

# Row 34 - LidoDAOToken
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.541"
$ws.Range("E34").Value = "  +0.91%  "

# Row 35 - HuobiToken
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.371"
$ws.Range("E35").Value = "  +0.31%  "

# Row 36 - ARBITRUM
$ws.Range("E36").Value = "  +0.67%  "

# Row 37 - ImmutableX
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.5512"
$ws.Range("E37").Value = "  +1.85%  "

# Row 38 - MXToken
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.538"
$ws.Range("E38").Value = "  -1.46%  "

# Row 39 - Maker
$ws.Range("D39").Value = "1.116.30"
$ws.Range("E39").Value = "  -2.73%  "

# Row 40 - VeChain
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.01556"
$ws.Range("E40").Value = "  +0.25%  "

# Row 41 - PaxDollar
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.001"
$ws.Range("E41").Value = "  +0.14%  "

# Row 42 - FraxShare
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.597"
$ws.Range("E42").Value = "  +2.73%  "

# Row 43 - TrustWalletToken
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.7954"
$ws.Range("E43").Value = "  -1.82%  "

# Row 44 - Quant
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "97.55"
$ws.Range("E44").Value = "  -1.29%  "

# Row 45 - RocketPoolETH
$ws.Range("D45").Value = "1.768.76"
$ws.Range("E45").Value = "  -0.41%  "

# Row 46 - BabyDogeCoin
$ws.Range("E46").Value = "  -7.29%  "

# Row 47 - Mantle
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.4446"
$ws.Range("E47").Value = "  -1.61%  "

# Row 48 - Frax
$ws.Range("E48").Value = "  +0.22%  "

# Row 49 - Aave
$ws.Range("E49").Value = "  -0.02%  "

# Row 50 - Cronos
$ws.Range("E50").Value = "  +1.60%  "

# Row 51 - EnergySwap
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "7.564"
$ws.Range("E51").Value = "  +2.84%  "
